$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [string][char]8323

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.500.73'
$ws.Range("E2").Value = '  +0.64%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.500.42'
$ws.Range("E3").Value = '  -0.42%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.71'
$ws.Range("E5").Value = '  +4.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.73'
$ws.Range("E6").Value = '  -2.67%  '

# Row 7
$ws.Range("E7").Value = '  -1.16%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.496.05'
$ws.Range("E8").Value = '  -0.34%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.03%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.198'
$ws.Range("E10").Value = '  +4.56%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.67'
$ws.Range("E11").Value = '  -0.53%  '

# Row 12
$ws.Range("E12").Value = '  -3.79%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.94'
$ws.Range("E13").Value = '  -0.77%  '

# Row 14
$ws.Range("E14").Value = '  +0.35%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.081.08'
$ws.Range("E15").Value = '  +0.00%  '

# Row 16
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.29'
$ws.Range("E16").Value = '  -6.15%  '

# Row 17
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '611.50'
$ws.Range("E17").Value = '  -8.82%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.502.17'
$ws.Range("E18").Value = '  -0.23%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.537.43'
$ws.Range("E19").Value = '  +0.72%  '

# Row 20
$ws.Range("E20").Value = '  -2.05%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.20'
$ws.Range("E21").Value = '  -2.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.17'
$ws.Range("E22").Value = '  -9.32%  '

# Row 23
$ws.Range("E23").Value = '  -2.94%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.79'
$ws.Range("E24").Value = '  -3.18%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.50'
$ws.Range("E25").Value = '  -3.01%  '

# Row 26
$ws.Range("E26").Value = '  -0.38%  '

# Row 27
$ws.Range("E27").Value = '  +0.03%  '

# Row 28
$ws.Range("E28").Value = '  -2.57%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.18'
$ws.Range("E29").Value = '  -2.78%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.08'
$ws.Range("E30").Value = '  +0.32%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.40'
$ws.Range("E31").Value = '  -4.11%  '

# Row 32
$ws.Range("E32").Value = '  -4.62%  '

# Row 33
$ws.Range("E33").Value = '  -2.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.89'
$ws.Range("E34").Value = '  -6.03%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '553.71'
$ws.Range("E35").Value = '  -4.18%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.73'
$ws.Range("E36").Value = '  -2.07%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.48'
$ws.Range("E37").Value = '  -3.51%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.72'
$ws.Range("E38").Value = '  -0.74%  '

# Row 39
$ws.Range("E39").Value = '  -4.43%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.09%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0447'
$ws.Range("E41").Value = '  +1.93%  '

# Row 42
$ws.Range("E42").Value = '  +1.35%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.330.07'
$ws.Range("E43").Value = '  -2.72%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.324'
$ws.Range("E44").Value = '  -3.94%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.91'
$ws.Range("E45").Value = '  -1.66%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0' + $sub3 + '0696'
$ws.Range("E46").Value = '  -1.50%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.90'
$ws.Range("E47").Value = '  -0.88%  '

# Row 48
$ws.Range("E48").Value = '  -0.08%  '

# Row 49
$ws.Range("E49").Value = '  -3.52%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '135.82'
$ws.Range("E50").Value = '  +3.65%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.67'
$ws.Range("E51").Value = '  +7.35%  '
